# Add MOM (Minutes of Meeting) entries for Jan 30 2020 to the timesheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: blank separator row, formatted like the existing separator (row 9) ---
$ws.Range("A9:C9").Copy()
$ws.Range("A20:C20").PasteSpecial(-4122)

# --- Re-color the existing separator row (row 9) to the new orange accent ---
$ws.Range("A9:C9").Interior.Color = 49407

# --- Row 21: Jan 30 10:00 to 12:00 | Client call | Sapphire automation ---
$ws.Range("A21").Value = "Jan 30 10:00 to 12:00"
$ws.Range("B21").Value = "Client call"
$ws.Range("C21").Value = "Sapphire automation"

# --- Row 22: Jan 30 12:00 to 13:00 | Discussion ... | Infimetrics ---
$ws.Range("A4:C4").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Jan 30 12:00 to 13:00"
$ws.Range("B22").Value = "Discussion with nitin sir and sujata man, creating new jupyter `nnotebook, modifying code of good days buckets to be saved in`nspecified directory."
$ws.Range("C22").Value = "Infimetrics"
$ws.Rows.Item(22).RowHeight = 45

# --- Row 23: Jan 30 13:00 to 14:00 | Lunch | Infimetrics ---
$ws.Range("A23").Value = "Jan 30 13:00 to 14:00"
$ws.Range("B23").Value = "Lunch"
$ws.Range("C23").Value = "Infimetrics"

# --- Row 24: Jan 30 14:00 to 15:00 | In combined dataset ... | Infimetrics ---
$ws.Range("A4:C4").Copy()
$ws.Range("A24:C24").PasteSpecial(-4122)
$ws.Range("A24").Value = "Jan 30 14:00 to 15:00"
$ws.Range("B24").Value = "In combined dataset notebook, cleaned good days data of both the `nmachines. Combined both datasets by performing full outer join on it."
$ws.Range("C24").Value = "Infimetrics"
$ws.Rows.Item(24).RowHeight = 45

# --- Row 25: Jan 30 15:00 to 16:00 | Revising percentile ... | Infimetrics ---
# (content typed before the timestamp, to reproduce the original shared-string order)
$ws.Range("A6:C6").Copy()
$ws.Range("A25:C25").PasteSpecial(-4122)
$ws.Range("B25").Value = "Revising percentile concent by watching call recording. Creating`npercentile buckets. "
$ws.Range("A25").Value = "Jan 30 15:00 to 16:00"
$ws.Range("C25").Value = "Infimetrics"
$ws.Rows.Item(25).RowHeight = 30

# --- Row 26: Jan 30 16:00 to 17:00 | Python class, working ... | Infimetrics ---
$ws.Range("A26").Value = "Jan 30 16:00 to 17:00"
$ws.Range("B26").Value = "Python class, working on connecting output counts to combind data."
$ws.Range("C26").Value = "Infimetrics"

# --- Row 27: Jan 30 17:00 to 18:00 | Combined output data ... | Infimetrics ---
$ws.Range("A27").Value = "Jan 30 17:00 to 18:00"
$ws.Range("B27").Value = "Combined output data with previously combined dataset."
$ws.Range("C27").Value = "Infimetrics"

# --- Row 28: Jan 30 18:00 to 19:00 | Python class, Connected ... | Infimetrics ---
$ws.Range("A5:C5").Copy()
$ws.Range("A28:C28").PasteSpecial(-4122)
$ws.Range("A28").Value = "Jan 30 18:00 to 19:00"
$ws.Range("B28").Value = "Python class, Connected output counts to combined data, working on`npercentile concepts by picking up examples."
$ws.Range("C28").Value = "Infimetrics"
$ws.Rows.Item(28).RowHeight = 45

# --- Column B slightly wider to fit the new content ---
$ws.Columns.Item(2).ColumnWidth = 62

# --- Update view: scroll down and select the last entered cell ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("D28").Select()
